# "Generate Report for Archive" — refresh the localization status report:
# the "zh-cn"/"de-de" status moves from "Ready for handoff" to "In Translation"
# everywhere it is referenced (Overview summary sheet + per-locale detail sheets).
# Excel auto-shrinks the now-narrower "Status" columns to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# Re-fit the status columns now that the text is shorter.
$overview.Range("E1:F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
